# Se agrega actualizacion automatica de productos al insertar uno
#
# Simula la sincronizacion automatica del listado de inventario contra el
# sistema de productos: al insertar un nuevo producto, la hoja completa se
# vuelve a poblar con el listado actualizado (altas nuevas + reordenamiento
# que entrega el backend).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Listado completo y actualizado de productos (ID, Nombre, Cantidad, Precio, Unidad, Categoria)
$productos = @(
    @(2, 'Mesa comedor 6p', 5, 549.5, 'juego', 'Comedor'),
    @(3, 'Silla giratoria', 20, 129.99, 'unidad', 'Oficina'),
    @(5, 'Juego de terraza', 3, 899.9, 'juego', 'Exteriores'),
    @(6, 'Lámpara de pie', 15, 85.75, 'unidad', 'Sala'),
    @(7, 'Escritorio ejecutivo', 8, 399, 'unidad', 'Oficina'),
    @(8, 'Velador doble', 10, 120, 'par', 'Dormitorio'),
    @(12, 'Zapatero zeta', 2, 120, 'unidad', 'Dormitorio'),
    @(4, 'Cama Queen', 10, 699, 'unidad', 'Dormitorio'),
    @(11, 'Zapatera', 10, 250, 'juego', 'Sala'),
    @(1, 'Sofá modular', 10, 799.99, 'unidad', 'Sala'),
    @(9, 'Silla plegable', 21, 45.99, 'unidad', 'Exteriores'),
    @(10, 'Repisa flotante', 30, 29.95, 'pieza', 'Oficina'),
    @(13, '$%&', 2, 2, 'par', 'Comedor'),
    @(14, 'Ropero Estela', 2, 300, 'unidad', 'Dormitorio'),
    @(15, 'Ropero choco', 23, 235, 'unidad', 'Dormitorio'),
    @(16, 'Comoda Lunaria', 23, 230, 'unidad', 'Dormitorio'),
    @(17, 'PRUEBA123', 212, 23, 'par', 'Comedor'),
    @(18, 'prueba4', 23, 23, 'unidad', 'Sala'),
    @(19, 'Prueba5', 1232, 231, 'juego', 'Oficina'),
    @(21, 'hola', 23, 234, 'par', 'Oficina'),
    @(20, 'hola123', 123, 123, 'juego', 'Oficina')
)

# Limpia el rango de datos previo (debajo del encabezado) antes de reescribirlo
$filasPrevias = $ws.UsedRange.Rows.Count
if ($filasPrevias -ge 2) {
    $ws.Range("A2:F" + $filasPrevias).ClearContents()
}

# Vuelca el listado actualizado de productos, fila a fila, a partir de la fila 2
for ($i = 0; $i -lt $productos.Length; $i++) {
    $fila = $productos[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $fila[0]
    $ws.Cells.Item($r, 2).Value2 = $fila[1]
    $ws.Cells.Item($r, 3).Value2 = $fila[2]
    $ws.Cells.Item($r, 4).Value2 = $fila[3]
    $ws.Cells.Item($r, 5).Value2 = $fila[4]
    $ws.Cells.Item($r, 6).Value2 = $fila[5]
}
